# DIV-4872 - Add UserRole column to CaseTypeTab
#
# Inserts a new "UserRole" column into the CaseTypeTab sheet, between the
# existing CaseFieldID column (H) and the TabFieldDisplayOrder column
# (previously I, now shifted right to J). The new column gets a header
# comment (row 2) explaining the MaxLength / role-restriction rules, and
# a sample column name (row 3) of "UserRole".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CaseTypeTab")

# CaseTypeTab becomes the active sheet/tab (previously FixedLists was
# selected); switching to it here drives that.
$ws.Activate()

# Insert a new column at I; Excel's default column-insert formatting is
# inherited from the column to the left, which matches the styles already
# used by the adjacent CaseFieldID column (H) for the shifted-in cells.
$ws.Columns("I:I").Insert()

# Row 2 (column header / help text) - new comment for the UserRole column.
$ws.Range("I2").Value = "MaxLength: 100. No entry for role means no role restriction for that tab. Enter role on a single row per tab"

# Row 3 (field name row) - the new column's name.
$ws.Range("I3").Value = "UserRole"

# Match the updated selection on this sheet.
$ws.Range("I3").Select()
